$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/02_pitito3.wav"
$ws.Range("B2").Value = "pngimages/02_pallet.png"

$ws.Range("A3").Value = "trainingaudio/03_kikita3.wav"
$ws.Range("B3").Value = "pngimages/03_box.png"

$ws.Range("A4").Value = "trainingaudio/24_takopa1.wav"
$ws.Range("B4").Value = "pngimages/24_banana.png"

$ws.Range("A5").Value = "trainingaudio/15_kopota3.wav"
$ws.Range("B5").Value = "pngimages/15_barrel.png"

$ws.Range("A6").Value = "trainingaudio/20_tatito1.wav"
$ws.Range("B6").Value = "pngimages/20_pizza.png"

$ws.Range("A7").Value = "trainingaudio/13_kopopi1.wav"
$ws.Range("B7").Value = "pngimages/13_toast.png"
